# edit.ps1 -- apply the "Updated cryptos list" refresh to Sheet1
# Re-applies new Price (D) / Volume(1h) (E) values scraped for this run,
# and two pairs of rows got re-ranked (Chainlink/Polkadot, PolygonEcosystemToken/
# EthereumClassic, Aave/WhiteBITCoin swapped order) so Coin (B) / Link (C) also change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain decimal-looking prices (e.g. "2.41") would be auto-parsed as numbers by
# Excel if written straight to .Value, which would also silently stamp a new
# number-format style onto the cell. Force them through as literal text by
# flipping to a Text format for the write and then clearing the style again so
# the cell keeps its original (unstyled) look -- matches the source being
# t="inlineStr" with no "s" attribute throughout.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Formula = "95.600.93"
$ws.Range("E2").Formula = "  -1.70%  "

# Row 3
$ws.Range("D3").Formula = "3.620.24"
$ws.Range("E3").Formula = "  -2.39%  "

# Row 4
Set-TextValue $ws.Range("D4") "2.41"
$ws.Range("E4").Formula = "  +26.24%  "

# Row 5
$ws.Range("E5").Formula = "  +0.04%  "

# Row 6
Set-TextValue $ws.Range("D6") "224.90"
$ws.Range("E6").Formula = "  -5.71%  "

# Row 7
Set-TextValue $ws.Range("D7") "637.91"
$ws.Range("E7").Formula = "  -3.29%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.418"
$ws.Range("E8").Formula = "  -1.59%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.09"
$ws.Range("E9").Formula = "  +2.51%  "

# Row 10
$ws.Range("E10").Formula = "  +0.06%  "

# Row 11
$ws.Range("D11").Formula = "3.615.29"
$ws.Range("E11").Formula = "  -2.48%  "

# Row 12
Set-TextValue $ws.Range("D12") "48.14"
$ws.Range("E12").Formula = "  +8.28%  "

# Row 13
$ws.Range("E13").Formula = "  -0.17%  "

# Row 14
$ws.Range("E14").Formula = "  -9.53%  "

# Row 15
Set-TextValue $ws.Range("D15") "6.54"
$ws.Range("E15").Formula = "  -4.50%  "

# Row 16
$ws.Range("D16").Formula = "4.296.78"
$ws.Range("E16").Formula = "  -2.32%  "

# Row 17
$ws.Range("D17").Formula = "95.386.79"
$ws.Range("E17").Formula = "  -1.72%  "

# Row 18
$ws.Range("B18").Formula = "Polkadot"
$ws.Range("C18").Formula = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D18") "8.80"
$ws.Range("E18").Formula = "  -3.40%  "

# Row 19
$ws.Range("B19").Formula = "Chainlink"
$ws.Range("C19").Formula = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D19") "20.86"
$ws.Range("E19").Formula = "  +11.80%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.00"
$ws.Range("E20").Formula = "  +7.87%  "

# Row 21
$ws.Range("D21").Formula = "3.615.17"
$ws.Range("E21").Formula = "  -2.27%  "

# Row 22
$ws.Range("E22").Formula = "  +2.70%  "

# Row 23
Set-TextValue $ws.Range("D23") "507.19"
$ws.Range("E23").Formula = "  -2.64%  "

# Row 24
$ws.Range("E24").Formula = "  -5.57%  "

# Row 25
Set-TextValue $ws.Range("D25") "0.247"
$ws.Range("E25").Formula = "  +26.35%  "

# Row 26
Set-TextValue $ws.Range("D26") "121.69"
$ws.Range("E26").Formula = "  +18.98%  "

# Row 27
$ws.Range("E27").Formula = "  -8.37%  "

# Row 28
Set-TextValue $ws.Range("D28") "6.77"
$ws.Range("E28").Formula = "  -1.98%  "

# Row 29
$ws.Range("D29").Formula = "3.819.35"
$ws.Range("E29").Formula = "  -2.25%  "

# Row 30
Set-TextValue $ws.Range("D30") "12.73"
$ws.Range("E30").Formula = "  -6.05%  "

# Row 31
Set-TextValue $ws.Range("D31") "13.02"
$ws.Range("E31").Formula = "  +1.53%  "

# Row 32
Set-TextValue $ws.Range("D32") "2.99"
$ws.Range("E32").Formula = "  -1.59%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.181"
$ws.Range("E34").Formula = "  -4.69%  "

# Row 35
$ws.Range("B35").Formula = "EthereumClassic"
$ws.Range("C35").Formula = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D35") "32.78"
$ws.Range("E35").Formula = "  +1.54%  "

# Row 36
$ws.Range("B36").Formula = "PolygonEcosystemToken"
$ws.Range("C36").Formula = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D36") "0.613"
$ws.Range("E36").Formula = "  +3.05%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Formula = "  -0.15%  "

# Row 38
Set-TextValue $ws.Range("D38") "1.76"
$ws.Range("E38").Formula = "  -5.62%  "

# Row 39
Set-TextValue $ws.Range("D39") "44.19"
$ws.Range("E39").Formula = "  +10.15%  "

# Row 41
Set-TextValue $ws.Range("D41") "593.05"
$ws.Range("E41").Formula = "  -9.71%  "

# Row 42
Set-TextValue $ws.Range("D42") "8.35"
$ws.Range("E42").Formula = "  -5.83%  "

# Row 43
Set-TextValue $ws.Range("D43") "7.01"
$ws.Range("E43").Formula = "  +2.64%  "

# Row 44
$ws.Range("E44").Formula = "  -5.47%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.478"
$ws.Range("E45").Formula = "  -1.47%  "

# Row 46
Set-TextValue $ws.Range("D46") "0.0485"
$ws.Range("E46").Formula = "  +5.69%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.94"
$ws.Range("E47").Formula = "  -4.93%  "

# Row 48
$ws.Range("E48").Formula = "  -2.33%  "

# Row 49
Set-TextValue $ws.Range("D49") "2.28"
$ws.Range("E49").Formula = "  -2.00%  "

# Row 50
$ws.Range("B50").Formula = "WhiteBITCoin"
$ws.Range("C50").Formula = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D50") "23.50"
$ws.Range("E50").Formula = "  -0.54%  "

# Row 51
$ws.Range("B51").Formula = "Aave"
$ws.Range("C51").Formula = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D51") "223.38"
$ws.Range("E51").Formula = "  +8.79%  "
